$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet / tab ---------------------------------------------------
$ws.Name = "format_petugas"

# --- Insert a new column for "kode_kegiatan" (E), shifting old E..H -> F..I
$ws.Columns("E:E").Insert()

# --- Header row ------------------------------------------------------------
$ws.Range("E1").Value = "kode_kegiatan"

# --- Update existing data rows (2-5) for the new email addresses / kode ---
$ws.Range("A2").Value = "petugas13@gmail.com"
$ws.Range("E2").Value = "sk04"

$ws.Range("A3").Value = "petugas14@gmail.com"
$ws.Range("E3").Value = "sk05"

$ws.Range("A4").Value = "petugas15@gmail.com"
$ws.Range("E4").Value = "sk06"

$ws.Range("A5").Value = "petugas16@gmail.com"
$ws.Range("E5").Value = "sk05"

# --- New row 6 ---------------------------------------------------------
$ws.Range("A6").Value = "petugas17@gmail.com"
$ws.Range("B6").Value = "password"
$ws.Range("C6").Value = "id15"
$ws.Range("D6").Value = "petugas15"
$ws.Range("E6").Value = "sk06"
$ws.Range("F6").Value = "laki-laki"
$ws.Range("G6").Value = 38
$ws.Range("I6").Value = "Jl.Helvetia"

# --- no_wa column (now H) becomes numeric with a "0" number format --------
$ws.Range("H2:H6").NumberFormat = "0"
$ws.Range("H2").Value = 853900821201
$ws.Range("H3").Value = 853900821202
$ws.Range("H4").Value = 853900821203
$ws.Range("H5").Value = 853900821204
$ws.Range("H6").Value = 853900821205

# --- Highlighted style for the "password" column (B2:B6) ------------------
$pwd = $ws.Range("B2:B6")
$pwd.Font.Color = 5732356
$pwd.Font.Name = "Arial"
$pwd.Interior.Color = 15071953
$pwd.Borders.Weight = -4138
$pwd.Borders.Color = 13421772
$pwd.WrapText = $true

# --- Row heights ------------------------------------------------------------
$ws.Rows("1:1").RowHeight = 15.75
$ws.Rows("2:6").RowHeight = 30

# --- Column width for no_wa (H) --------------------------------------------
$ws.Columns("H:H").ColumnWidth = 13.140625

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection / view ------------------------------------------------------
$null = $ws.Range("B2:B6").Select()
